$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B188:B191").ClearContents()
$ws.Range("B234:B263").ClearContents()
$ws.Range("B312:B315").ClearContents()
$ws.Range("B518").ClearContents()
$ws.Range("B1052:B1053").ClearContents()
$ws.Range("B1938").ClearContents()
$ws.Range("B1976:B1977").ClearContents()
$ws.Range("B2029").ClearContents()
$ws.Range("B2032").ClearContents()
$ws.Range("B4239:B4242").ClearContents()
$ws.Range("B4771:B4772").ClearContents()
$ws.Range("B4820:B4864").ClearContents()
$ws.Range("B6448:B6449").ClearContents()
$ws.Range("B6566:B6568").ClearContents()
$ws.Range("B7414:B7416").ClearContents()
$ws.Range("B8345:B8485").ClearContents()
Write-Output "cleared"
